$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7-10 down to 8-11
$ws.Rows("7:7").Insert()

# Populate the new row 7 with the new weekly price record
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = 44894
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100103
$ws.Range("H7").Value = "Frutos de hueso (carozo)"
$ws.Range("I7").Value = 100103003
$ws.Range("J7").Value = "Damasco"
$ws.Range("K7").Value = "Castle Brite"
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 130
$ws.Range("N7").Value = 19000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 19462
$ws.Range("Q7").Value = "$/caja 16 kilos granel"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1216
$ws.Range("T7").Value = 16
